$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.278"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05685"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.322"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8067"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8939"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1429"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'9WazirXWRX"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.07454"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10MandalaExchangeTokenMDX"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'ProBitToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.1370"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11ProBitTokenPROB"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03063"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Value = "'0.09395"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.879"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001585"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04803"
$ws.Range("D18").Style = "Normal"
$ws.Range("B20").Value = "'One"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.01097"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'19OneONE"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'TigerCash"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'0.006418"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'20TigerCashTCH"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'HotbitToken"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.004972"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'21HotbitTokenHTB"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'BitKan"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.0009972"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'22BitKanKAN"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'NitroEx"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.0001501"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'23NitroExNTX"
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'LEO"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'3.686"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'24LEOLEO"
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'BTSEToken"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'2.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'25BTSETokenBTSE"
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'BitpandaEcosystemToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'0.3257"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'26BitpandaEcosystemTokenBEST"
$ws.Range("E27").Style = "Normal"
$ws.Range("D41").Value = "'0.006845"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1067"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002819"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008755"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005594"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.4994"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.2027"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
